$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86 (hunk 0)
$ws.Cells.Item(86, 8).Value = 3763958.2  # H86: 4053316.5 -> 3763958.2
$ws.Cells.Item(86, 9).Value = 3926.8  # I86: 4333.75 -> 3926.8
$ws.Cells.Item(86, 11).Value = 3926.8  # K86: 4333.75 -> 3926.8
$ws.Cells.Item(86, 13).Value = -2803.8  # M86: -3210.75 -> -2803.8

# Row 89 (hunk 1)
$ws.Cells.Item(89, 8).Value = 3763958.2  # H89: 4053316.5 -> 3763958.2
$ws.Cells.Item(89, 9).Value = 3926.8  # I89: 4333.75 -> 3926.8
$ws.Cells.Item(89, 11).Value = 19634  # K89: 21668.75 -> 19634
$ws.Cells.Item(89, 13).Value = -14018  # M89: -16052.75 -> -14018

# Row 100 (hunk 2)
$ws.Cells.Item(100, 8).Value = 2882.1538  # H100: 2803.3572 -> 2882.1538
$ws.Cells.Item(100, 9).Value = 3151.3333  # I100: 2955.2856 -> 3151.3333
$ws.Cells.Item(100, 11).Value = 3151.3333  # K100: 2955.2856 -> 3151.3333
$ws.Cells.Item(100, 13).Value = -2610.3333  # M100: -2414.2856 -> -2610.3333

# Row 135 (hunk 3)
$ws.Cells.Item(135, 8).Value = 2250.7144  # H135: 1156.4 -> 2250.7144
$ws.Cells.Item(135, 9).Value = 959.1667  # I135: 661.7 -> 959.1667
$ws.Cells.Item(135, 10).Value = 10000  # J135: 2145.8 -> 10000
$ws.Cells.Item(135, 11).Value = 8632.5003  # K135: 5955.3 -> 8632.5003
$ws.Cells.Item(135, 12).Value = 90000  # L135: 19312.2 -> 90000
$ws.Cells.Item(135, 13).Value = -6097.5003  # M135: -3420.3 -> -6097.5003
$ws.Cells.Item(135, 14).Value = -95070  # N135: -24382.2 -> -95070

# Row 138 (hunk 4)
$ws.Cells.Item(138, 8).Value = 7121.625  # H138: 7102.706 -> 7121.625
$ws.Cells.Item(138, 10).Value = 8120.636  # J138: 8010.5835 -> 8120.636
$ws.Cells.Item(138, 12).Value = 24361.908  # L138: 24031.7505 -> 24361.908
$ws.Cells.Item(138, 14).Value = -34641.908  # N138: -34311.75049999999 -> -34641.908

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 5)
$ws.Cells.Item(32, 8).Value = 2090.2222  # H32: 2174.0566 -> 2090.2222
$ws.Cells.Item(32, 9).Value = 1949.64  # I32: 2037.449 -> 1949.64
$ws.Cells.Item(32, 11).Value = 1949.64  # K32: 2037.449 -> 1949.64
$ws.Cells.Item(32, 13).Value = -1662.64  # M32: -1750.449 -> -1662.64

# Row 61 (hunk 6)
$ws.Cells.Item(61, 8).Value = 4223.4443  # H61: 4251.222 -> 4223.4443
$ws.Cells.Item(61, 9).Value = 3404.2  # I61: 3671.3333 -> 3404.2
$ws.Cells.Item(61, 10).Value = 5247.5  # J61: 4831.1113 -> 5247.5
$ws.Cells.Item(61, 11).Value = 3404.2  # K61: 3671.3333 -> 3404.2
$ws.Cells.Item(61, 12).Value = 5247.5  # L61: 4831.1113 -> 5247.5
$ws.Cells.Item(61, 13).Value = -3192.2  # M61: -3459.3333 -> -3192.2
$ws.Cells.Item(61, 14).Value = -5671.5  # N61: -5255.1113 -> -5671.5

# Row 97 (hunk 7)
$ws.Cells.Item(97, 8).Value = 2588.647  # H97: 2095.652 -> 2588.647
$ws.Cells.Item(97, 9).Value = 2935.5  # I97: 2264.5 -> 2935.5
$ws.Cells.Item(97, 11).Value = 2935.5  # K97: 2264.5 -> 2935.5
$ws.Cells.Item(97, 13).Value = -2439.5  # M97: -1768.5 -> -2439.5

# Row 102 (hunk 8)
$ws.Cells.Item(102, 8).Value = 4751.75  # H102: 4999.5 -> 4751.75
$ws.Cells.Item(102, 9).Value = 4669  # I102: 4999 -> 4669
$ws.Cells.Item(102, 11).Value = 4669  # K102: 4999 -> 4669
$ws.Cells.Item(102, 13).Value = -3047  # M102: -3377 -> -3047

# Row 122 (hunk 9)
$ws.Cells.Item(122, 8).Value = 38465130  # H122: 38465136 -> 38465130
$ws.Cells.Item(122, 9).Value = 200002640  # I122: 250003060 -> 200002640
$ws.Cells.Item(122, 10).Value = 3816.0476  # J122: 3697.0908 -> 3816.0476
$ws.Cells.Item(122, 11).Value = 600007920  # K122: 750009180 -> 600007920
$ws.Cells.Item(122, 12).Value = 11448.1428  # L122: 11091.2724 -> 11448.1428
$ws.Cells.Item(122, 13).Value = -600005470  # M122: -750006730 -> -600005470
$ws.Cells.Item(122, 14).Value = -16348.1428  # N122: -15991.2724 -> -16348.1428

# Row 136 (hunk 10)
$ws.Cells.Item(136, 8).Value = 4223.4443  # H136: 4251.222 -> 4223.4443
$ws.Cells.Item(136, 9).Value = 3404.2  # I136: 3671.3333 -> 3404.2
$ws.Cells.Item(136, 10).Value = 5247.5  # J136: 4831.1113 -> 5247.5
$ws.Cells.Item(136, 11).Value = 10212.6  # K136: 11013.9999 -> 10212.6
$ws.Cells.Item(136, 12).Value = 15742.5  # L136: 14493.3339 -> 15742.5
$ws.Cells.Item(136, 13).Value = -7662.599999999999  # M136: -8463.999899999999 -> -7662.599999999999
$ws.Cells.Item(136, 14).Value = -20842.5  # N136: -19593.3339 -> -20842.5

# Row 139 (hunk 11)
$ws.Cells.Item(139, 8).Value = 50000  # H139: 49612.5 -> 50000
$ws.Cells.Item(139, 10).Value = 50000  # J139: 49612.5 -> 50000
$ws.Cells.Item(139, 12).Value = 50000  # L139: 49612.5 -> 50000
$ws.Cells.Item(139, 14).Value = -60280  # N139: -59892.5 -> -60280

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (hunk 12)
$ws.Cells.Item(20, 8).Value = 3321.5625  # H20: 3411.6428 -> 3321.5625
$ws.Cells.Item(20, 9).Value = 3171  # I20: 3277.6667 -> 3171
$ws.Cells.Item(20, 11).Value = 3171  # K20: 3277.6667 -> 3171
$ws.Cells.Item(20, 13).Value = -2924  # M20: -3030.6667 -> -2924

# Row 86 (hunk 13)
$ws.Cells.Item(86, 8).Value = 1640654.1  # H86: 1804269.5 -> 1640654.1
$ws.Cells.Item(86, 10).Value = 5216.3335  # J86: 5359.6 -> 5216.3335
$ws.Cells.Item(86, 12).Value = 5216.3335  # L86: 5359.6 -> 5216.3335
$ws.Cells.Item(86, 14).Value = -7462.3335  # N86: -7605.6 -> -7462.3335

# Row 89 (hunk 14)
$ws.Cells.Item(89, 8).Value = 1640654.1  # H89: 1804269.5 -> 1640654.1
$ws.Cells.Item(89, 10).Value = 5216.3335  # J89: 5359.6 -> 5216.3335
$ws.Cells.Item(89, 12).Value = 26081.6675  # L89: 26798 -> 26081.6675
$ws.Cells.Item(89, 14).Value = -37313.6675  # N89: -38030 -> -37313.6675

# Row 105 (hunk 15)
$ws.Cells.Item(105, 8).Value = 15154093  # H105: 18521682 -> 15154093
$ws.Cells.Item(105, 9).Value = 2336.5715  # I105: 3114.75 -> 2336.5715
$ws.Cells.Item(105, 10).Value = 41669668  # J105: 33336534 -> 41669668
$ws.Cells.Item(105, 11).Value = 2336.5715  # K105: 3114.75 -> 2336.5715
$ws.Cells.Item(105, 12).Value = 41669668  # L105: 33336534 -> 41669668
$ws.Cells.Item(105, 13).Value = -589.5715  # M105: -1367.75 -> -589.5715
$ws.Cells.Item(105, 14).Value = -41673162  # N105: -33340028 -> -41673162

# Row 107 (hunk 16)
$ws.Cells.Item(107, 8).Value = 2984.2  # H107: 4370.3335 -> 2984.2
$ws.Cells.Item(107, 9).Value = 2984.2  # I107: 4370.3335 -> 2984.2
$ws.Cells.Item(107, 11).Value = 2984.2  # K107: 4370.3335 -> 2984.2
$ws.Cells.Item(107, 13).Value = -1064.2  # M107: -2450.3335 -> -1064.2

# Row 134 (hunk 17)
$ws.Cells.Item(134, 8).Value = 82889.62  # H134: 82965.53999999999 -> 82889.62
$ws.Cells.Item(134, 9).Value = 5414  # I134: 5660.75 -> 5414
$ws.Cells.Item(134, 11).Value = 16242  # K134: 16982.25 -> 16242
$ws.Cells.Item(134, 13).Value = -13707  # M134: -14447.25 -> -13707

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (hunk 18)
$ws.Cells.Item(7, 8).Value = 406.5862  # H7: 364.30304 -> 406.5862
$ws.Cells.Item(7, 9).Value = 434.26923  # I7: 394.2069 -> 434.26923
$ws.Cells.Item(7, 10).Value = 166.66667  # J7: 147.5 -> 166.66667
$ws.Cells.Item(7, 11).Value = 434.26923  # K7: 394.2069 -> 434.26923
$ws.Cells.Item(7, 12).Value = 166.66667  # L7: 147.5 -> 166.66667
$ws.Cells.Item(7, 13).Value = -321.26923  # M7: -281.2069 -> -321.26923
$ws.Cells.Item(7, 14).Value = -392.66667  # N7: -373.5 -> -392.66667

# Row 31 (hunk 19)
$ws.Cells.Item(31, 8).Value = 2124170.5  # H31: 1133397.2 -> 2124170.5
$ws.Cells.Item(31, 9).Value = 8475182  # I31: 5650420.5 -> 8475182
$ws.Cells.Item(31, 10).Value = 7166.5  # J31: 4141.4165 -> 7166.5
$ws.Cells.Item(31, 11).Value = 8475182  # K31: 5650420.5 -> 8475182
$ws.Cells.Item(31, 12).Value = 7166.5  # L31: 4141.4165 -> 7166.5
$ws.Cells.Item(31, 13).Value = -8474887  # M31: -5650125.5 -> -8474887
$ws.Cells.Item(31, 14).Value = -7756.5  # N31: -4731.4165 -> -7756.5

# Row 34 (hunk 20)
$ws.Cells.Item(34, 8).Value = 2124170.5  # H34: 1133397.2 -> 2124170.5
$ws.Cells.Item(34, 9).Value = 8475182  # I34: 5650420.5 -> 8475182
$ws.Cells.Item(34, 10).Value = 7166.5  # J34: 4141.4165 -> 7166.5
$ws.Cells.Item(34, 11).Value = 8475182  # K34: 5650420.5 -> 8475182
$ws.Cells.Item(34, 12).Value = 7166.5  # L34: 4141.4165 -> 7166.5
$ws.Cells.Item(34, 13).Value = -8474980  # M34: -5650218.5 -> -8474980
$ws.Cells.Item(34, 14).Value = -7570.5  # N34: -4545.4165 -> -7570.5

# Row 58 (hunk 21)
$ws.Cells.Item(58, 8).Value = 3682.913  # H58: 3491.6538 -> 3682.913
$ws.Cells.Item(58, 9).Value = 1729.1428  # I58: 1807.3846 -> 1729.1428
$ws.Cells.Item(58, 10).Value = 6722.1113  # J58: 5175.923 -> 6722.1113
$ws.Cells.Item(58, 11).Value = 1729.1428  # K58: 1807.3846 -> 1729.1428
$ws.Cells.Item(58, 12).Value = 6722.1113  # L58: 5175.923 -> 6722.1113
$ws.Cells.Item(58, 13).Value = -1526.1428  # M58: -1604.3846 -> -1526.1428
$ws.Cells.Item(58, 14).Value = -7128.1113  # N58: -5581.923 -> -7128.1113

# Row 86 (hunk 22)
$ws.Cells.Item(86, 8).Value = 2796.1177  # H86: 2948.4707 -> 2796.1177
$ws.Cells.Item(86, 9).Value = 2389.75  # I86: 2515.8572 -> 2389.75
$ws.Cells.Item(86, 10).Value = 3157.3333  # J86: 3251.3 -> 3157.3333
$ws.Cells.Item(86, 11).Value = 2389.75  # K86: 2515.8572 -> 2389.75
$ws.Cells.Item(86, 12).Value = 3157.3333  # L86: 3251.3 -> 3157.3333
$ws.Cells.Item(86, 13).Value = -1266.75  # M86: -1392.8572 -> -1266.75
$ws.Cells.Item(86, 14).Value = -5403.3333  # N86: -5497.3 -> -5403.3333

# Row 89 (hunk 23)
$ws.Cells.Item(89, 8).Value = 2796.1177  # H89: 2948.4707 -> 2796.1177
$ws.Cells.Item(89, 9).Value = 2389.75  # I89: 2515.8572 -> 2389.75
$ws.Cells.Item(89, 10).Value = 3157.3333  # J89: 3251.3 -> 3157.3333
$ws.Cells.Item(89, 11).Value = 11948.75  # K89: 12579.286 -> 11948.75
$ws.Cells.Item(89, 12).Value = 15786.6665  # L89: 16256.5 -> 15786.6665
$ws.Cells.Item(89, 13).Value = -6332.75  # M89: -6963.286 -> -6332.75
$ws.Cells.Item(89, 14).Value = -27018.6665  # N89: -27488.5 -> -27018.6665

# Row 105 (hunk 24)
$ws.Cells.Item(105, 8).Value = 2851.5715  # H105: 2351.111 -> 2851.5715
$ws.Cells.Item(105, 9).Value = 3017  # I105: 2211.3333 -> 3017
$ws.Cells.Item(105, 10).Value = 2631  # J105: 2630.6667 -> 2631
$ws.Cells.Item(105, 11).Value = 3017  # K105: 2211.3333 -> 3017
$ws.Cells.Item(105, 12).Value = 2631  # L105: 2630.6667 -> 2631
$ws.Cells.Item(105, 13).Value = -1270  # M105: -464.3332999999998 -> -1270
$ws.Cells.Item(105, 14).Value = -6125  # N105: -6124.6667 -> -6125

# Row 132 (hunk 25)
$ws.Cells.Item(132, 8).Value = 3138.8635  # H132: 3264.8696 -> 3138.8635
$ws.Cells.Item(132, 9).Value = 2539.5293  # I132: 2541.7646 -> 2539.5293
$ws.Cells.Item(132, 10).Value = 5176.6  # J132: 5313.6665 -> 5176.6
$ws.Cells.Item(132, 11).Value = 7618.5879  # K132: 7625.293799999999 -> 7618.5879
$ws.Cells.Item(132, 12).Value = 15529.8  # L132: 15940.9995 -> 15529.8
$ws.Cells.Item(132, 13).Value = -5088.5879  # M132: -5095.293799999999 -> -5088.5879
$ws.Cells.Item(132, 14).Value = -20589.8  # N132: -21000.9995 -> -20589.8

# Row 136 (hunk 26)
$ws.Cells.Item(136, 8).Value = 3682.913  # H136: 3491.6538 -> 3682.913
$ws.Cells.Item(136, 9).Value = 1729.1428  # I136: 1807.3846 -> 1729.1428
$ws.Cells.Item(136, 10).Value = 6722.1113  # J136: 5175.923 -> 6722.1113
$ws.Cells.Item(136, 11).Value = 5187.428400000001  # K136: 5422.1538 -> 5187.428400000001
$ws.Cells.Item(136, 12).Value = 20166.3339  # L136: 15527.769 -> 20166.3339
$ws.Cells.Item(136, 13).Value = -2637.428400000001  # M136: -2872.1538 -> -2637.428400000001
$ws.Cells.Item(136, 14).Value = -25266.3339  # N136: -20627.769 -> -25266.3339

$ws = $wb.Worksheets.Item("CUL")
# Row 26 (hunk 27)
$ws.Cells.Item(26, 8).Value = 57358.6  # H26: 60801.938 -> 57358.6
$ws.Cells.Item(26, 9).Value = 250079.12  # I26: 250079.5 -> 250079.12
$ws.Cells.Item(26, 10).Value = 256.22223  # J26: 233.12 -> 256.22223
$ws.Cells.Item(26, 11).Value = 750237.36  # K26: 750238.5 -> 750237.36
$ws.Cells.Item(26, 12).Value = 768.66669  # L26: 699.36 -> 768.66669
$ws.Cells.Item(26, 13).Value = -749949.36  # M26: -749950.5 -> -749949.36
$ws.Cells.Item(26, 14).Value = -1344.66669  # N26: -1275.36 -> -1344.66669

# Row 34 (hunk 28)
$ws.Cells.Item(34, 8).Value = 7583.95  # H34: 51735.41 -> 7583.95
$ws.Cells.Item(34, 9).Value = 3029.6  # I34: 2378.2856 -> 3029.6
$ws.Cells.Item(34, 10).Value = 9102.066000000001  # J34: 74768.734 -> 9102.066000000001
$ws.Cells.Item(34, 11).Value = 9088.799999999999  # K34: 7134.8568 -> 9088.799999999999
$ws.Cells.Item(34, 12).Value = 27306.198  # L34: 224306.202 -> 27306.198
$ws.Cells.Item(34, 13).Value = -9004.799999999999  # M34: -7050.8568 -> -9004.799999999999
$ws.Cells.Item(34, 14).Value = -27474.198  # N34: -224474.202 -> -27474.198

# Row 68 (hunk 29)
$ws.Cells.Item(68, 8).Value = 2055.2708  # H68: 228250.48 -> 2055.2708
$ws.Cells.Item(68, 10).Value = 2069.0244  # J68: 263448.94 -> 2069.0244
$ws.Cells.Item(68, 12).Value = 6207.073199999999  # L68: 790346.8200000001 -> 6207.073199999999
$ws.Cells.Item(68, 14).Value = -7829.073199999999  # N68: -791968.8200000001 -> -7829.073199999999

# Row 71 (hunk 30)
$ws.Cells.Item(71, 8).Value = 2055.2708  # H71: 228250.48 -> 2055.2708
$ws.Cells.Item(71, 10).Value = 2069.0244  # J71: 263448.94 -> 2069.0244
$ws.Cells.Item(71, 12).Value = 18621.2196  # L71: 2371040.46 -> 18621.2196
$ws.Cells.Item(71, 14).Value = -26733.2196  # N71: -2379152.46 -> -26733.2196

# Row 107 (hunk 31)
$ws.Cells.Item(107, 8).Value = 26861.547  # H107: 26289.162 -> 26861.547
$ws.Cells.Item(107, 10).Value = 33014.35  # J107: 32135.344 -> 33014.35
$ws.Cells.Item(107, 12).Value = 99043.04999999999  # L107: 96406.03200000001 -> 99043.04999999999
$ws.Cells.Item(107, 14).Value = -102883.05  # N107: -100246.032 -> -102883.05

# Row 131 (hunk 32)
$ws.Cells.Item(131, 8).Value = 63759.758  # H131: 61902.117 -> 63759.758
$ws.Cells.Item(131, 9).Value = 72581.57000000001  # I131: 78118.62 -> 72581.57000000001
$ws.Cells.Item(131, 10).Value = 57259.473  # J131: 51863.332 -> 57259.473
$ws.Cells.Item(131, 11).Value = 217744.71  # K131: 234355.86 -> 217744.71
$ws.Cells.Item(131, 12).Value = 171778.419  # L131: 155589.996 -> 171778.419
$ws.Cells.Item(131, 13).Value = -212704.71  # M131: -229315.86 -> -212704.71
$ws.Cells.Item(131, 14).Value = -181858.419  # N131: -165669.996 -> -181858.419

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 33)
$ws.Cells.Item(2, 8).Value = 157.79411  # H2: 161.69698 -> 157.79411
$ws.Cells.Item(2, 10).Value = 221.2  # J2: 242.55556 -> 221.2
$ws.Cells.Item(2, 12).Value = 221.2  # L2: 242.55556 -> 221.2
$ws.Cells.Item(2, 14).Value = -447.2  # N2: -468.55556 -> -447.2

# Row 97 (hunk 34)
$ws.Cells.Item(97, 8).Value = 1538.1  # H97: 1384.68 -> 1538.1
$ws.Cells.Item(97, 9).Value = 1209.4375  # I97: 1145.35 -> 1209.4375
$ws.Cells.Item(97, 10).Value = 2852.75  # J97: 2342 -> 2852.75
$ws.Cells.Item(97, 11).Value = 1209.4375  # K97: 1145.35 -> 1209.4375
$ws.Cells.Item(97, 12).Value = 2852.75  # L97: 2342 -> 2852.75
$ws.Cells.Item(97, 13).Value = -713.4375  # M97: -649.3499999999999 -> -713.4375
$ws.Cells.Item(97, 14).Value = -3844.75  # N97: -3334 -> -3844.75

# Row 122 (hunk 35)
$ws.Cells.Item(122, 8).Value = 414873.4  # H122: 400124.3 -> 414873.4
$ws.Cells.Item(122, 9).Value = 619254.6  # I122: 586762.2 -> 619254.6
$ws.Cells.Item(122, 11).Value = 1857763.8  # K122: 1760286.6 -> 1857763.8
$ws.Cells.Item(122, 13).Value = -1855313.8  # M122: -1757836.6 -> -1855313.8

$ws = $wb.Worksheets.Item("LTW")
# Row 25 (hunk 36)
$ws.Cells.Item(25, 8).Value = 7529.125  # H25: 9106.6 -> 7529.125
$ws.Cells.Item(25, 10).Value = 6046.8  # J25: 7767 -> 6046.8
$ws.Cells.Item(25, 12).Value = 6046.8  # L25: 7767 -> 6046.8
$ws.Cells.Item(25, 14).Value = -6506.8  # N25: -8227 -> -6506.8

# Row 40 (hunk 37)
$ws.Cells.Item(40, 8).Value = 7502625  # H40: 6003100 -> 7502625
$ws.Cells.Item(40, 9).Value = 10001767  # I40: 8573086 -> 10001767
$ws.Cells.Item(40, 10).Value = 5200  # J40: 6466.6665 -> 5200
$ws.Cells.Item(40, 11).Value = 10001767  # K40: 8573086 -> 10001767
$ws.Cells.Item(40, 12).Value = 5200  # L40: 6466.6665 -> 5200
$ws.Cells.Item(40, 13).Value = -10001631  # M40: -8572950 -> -10001631
$ws.Cells.Item(40, 14).Value = -5472  # N40: -6738.6665 -> -5472

# Row 61 (hunk 38)
$ws.Cells.Item(61, 8).Value = 4279.579  # H61: 4461.8335 -> 4279.579
$ws.Cells.Item(61, 9).Value = 3427.577  # I61: 3524.72 -> 3427.577
$ws.Cells.Item(61, 10).Value = 6125.5835  # J61: 6591.636 -> 6125.5835
$ws.Cells.Item(61, 11).Value = 3427.577  # K61: 3524.72 -> 3427.577
$ws.Cells.Item(61, 12).Value = 6125.5835  # L61: 6591.636 -> 6125.5835
$ws.Cells.Item(61, 13).Value = -3225.577  # M61: -3322.72 -> -3225.577
$ws.Cells.Item(61, 14).Value = -6529.5835  # N61: -6995.636 -> -6529.5835

# Row 113 (hunk 39)
$ws.Cells.Item(113, 8).Value = 4279.579  # H113: 4461.8335 -> 4279.579
$ws.Cells.Item(113, 9).Value = 3427.577  # I113: 3524.72 -> 3427.577
$ws.Cells.Item(113, 10).Value = 6125.5835  # J113: 6591.636 -> 6125.5835
$ws.Cells.Item(113, 11).Value = 3427.577  # K113: 3524.72 -> 3427.577
$ws.Cells.Item(113, 12).Value = 6125.5835  # L113: 6591.636 -> 6125.5835
$ws.Cells.Item(113, 13).Value = -1257.577  # M113: -1354.72 -> -1257.577
$ws.Cells.Item(113, 14).Value = -10465.5835  # N113: -10931.636 -> -10465.5835

# Row 122 (hunk 40)
$ws.Cells.Item(122, 8).Value = 1231445.5  # H122: 558217.9 -> 1231445.5
$ws.Cells.Item(122, 9).Value = 0  # I122: 7334.7 -> 0
$ws.Cells.Item(122, 10).Value = 1231445.5  # J122: 1109101 -> 1231445.5
$ws.Cells.Item(122, 11).Value = 0  # K122: 22004.1 -> 0
$ws.Cells.Item(122, 12).Value = 3694336.5  # L122: 3327303 -> 3694336.5
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -19554.1
$ws.Cells.Item(122, 14).Value = -3699236.5  # N122: -3332203 -> -3699236.5

# Row 136 (hunk 41)
$ws.Cells.Item(136, 8).Value = 926871.7  # H136: 870807.75 -> 926871.7
$ws.Cells.Item(136, 9).Value = 1590557.4  # I136: 1431683.2 -> 1590557.4
$ws.Cells.Item(136, 11).Value = 4771672.199999999  # K136: 4295049.6 -> 4771672.199999999
$ws.Cells.Item(136, 13).Value = -4769122.199999999  # M136: -4292499.6 -> -4769122.199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (hunk 42)
$ws.Cells.Item(54, 8).Value = 101077  # H54: 75538.5 -> 101077
$ws.Cells.Item(54, 9).Value = 0  # I54: 50000 -> 0
$ws.Cells.Item(54, 11).Value = 0  # K54: 50000 -> 0
$ws.Cells.Item(54, 13).ClearContents()  # M54: was -49480

# Row 81 (hunk 43)
$ws.Cells.Item(81, 10).Value = 5000  # J81: 0 -> 5000
$ws.Cells.Item(81, 12).Value = 10000  # L81: 0 -> 10000
$ws.Cells.Item(81, 14).Value = -12122  # N81: None -> -12122

# Row 84 (hunk 44)
$ws.Cells.Item(84, 10).Value = 5000  # J84: 0 -> 5000
$ws.Cells.Item(84, 12).Value = 50000  # L84: 0 -> 50000
$ws.Cells.Item(84, 14).Value = -60608  # N84: None -> -60608

# Row 126 (hunk 45)
$ws.Cells.Item(126, 8).Value = 3487.111  # H126: 3585.5625 -> 3487.111
$ws.Cells.Item(126, 9).Value = 3058.25  # I126: 3130 -> 3058.25
$ws.Cells.Item(126, 11).Value = 9174.75  # K126: 9390 -> 9174.75
$ws.Cells.Item(126, 13).Value = -6704.75  # M126: -6920 -> -6704.75

# Row 136 (hunk 46)
$ws.Cells.Item(136, 8).Value = 505699.38  # H136: 418626.22 -> 505699.38
$ws.Cells.Item(136, 9).Value = 718251.1  # I136: 591617.7 -> 718251.1
$ws.Cells.Item(136, 10).Value = 208126.9  # J136: 173554.92 -> 208126.9
$ws.Cells.Item(136, 11).Value = 2154753.3  # K136: 1774853.1 -> 2154753.3
$ws.Cells.Item(136, 12).Value = 624380.7  # L136: 520664.76 -> 624380.7
$ws.Cells.Item(136, 13).Value = -2152203.3  # M136: -1772303.1 -> -2152203.3
$ws.Cells.Item(136, 14).Value = -629480.7  # N136: -525764.76 -> -629480.7
